$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column Z was an unused/empty column in the data rows (2-4), while the
#    header row already used it (no gap, header row only spans A-AH). Close
#    that gap by deleting the empty Z cells in rows 2-4 and shifting the
#    remaining data (Z..AI) left by one column, so the data lines back up
#    with the headers. (Do this before touching row 1, since the shift
#    operates on whole columns.)
$ws.Range("Z2:Z4").Delete(-4159)

# 2) Rename the header row (row 1) labels to the new uppercase/renamed values.
$headers = @(
    "FECHA SEGUIMIENTO", "FECHATERMINO", "ID", "SINIESTRO", "POLIZA", "AFECTADO",
    "TIPO DE CASO", "COBERTURA", "FECHA SINIESTRO", "ESTADO", "CIUDAD", "REGION",
    "UBICACION TALLER", "FINANCIADO", "REGIMEN FISCAL", "ESTATUS CLIENTE", "COMENTARIOS",
    "FECHA CARGA", "FECHA DECRETO", "USUARIO DE CARGA", "ESTATUS SEGUIMIENTO",
    "USUARIO ASIGNADO", "FECHA ASIGNACION", "MARCA", "TIPO", "MODELO", "NUMERO SERIE",
    "VALOR INDEMNIZACION", "VALOR COMERCIAL", "PLACAS", "ESTACION", " ESTATUS",
    "SUB ESTATUS", "USUARIO EN SEGUIMIENTO"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
